# Update the resampled/bootstrapped results in column A (rows 2-49)
# and slightly narrow column A's width, per the new run of the
# "generacion de algunos resultados" computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    2.0096818918040058,
    0.10419017175436893,
    0.06073459429499839,
    0.10419086645753386,
    0.46471381381104365,
    0.46471389022953691,
    -0.013365773974080092,
    -0.013365738697902568,
    0.23980343489296019,
    0.008111571558382916,
    -0.019554380766878279,
    0.0084885583754966114,
    0.085742461613731419,
    0.085742478840539899,
    -0.025265524575219321,
    -0.025265513780247363,
    -0.023781880370353222,
    0.026734490545885517,
    -0.023781433025612253,
    0.025784614172220267,
    0.0010160412959051776,
    0.0010160437866119018,
    -0.0266043839074874,
    -0.026604371092488845,
    -0.01501101430663007,
    -0.024666573486693193,
    -0.024641097271382584,
    -0.02466923616441908,
    -0.021941043027680091,
    -0.021941040828049856,
    -0.02381837558615885,
    -0.023818384041146626,
    -0.022560861058491109,
    -0.02265233644541105,
    -0.022691416489542315,
    -0.022675637466859262,
    -0.025175271490418535,
    -0.02517524525397458,
    -0.023766056087734674,
    -0.023766077738608553,
    -0.027372902024398661,
    -0.027372902212269269,
    -0.0264601265519449,
    -0.026460129125311784,
    -0.021872934979597754,
    -0.021872845896268551,
    -0.018632828500986312,
    -0.018632893238897167
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Narrow column A slightly (originally ~16.21875 character-units wide).
$ws.Columns.Item(1).ColumnWidth = 13.65
